# "Conditions use for different activity"
# - Expand "CreateOwn Test Data" (sheet1) with new columns G:O describing
#   activity-type / schedule / assignment / completion test data.
# - Insert a new blank-header "Sheet1" (just data row, no header) between
#   "CreateOwn Test Data" and "AddAchievement Test Data".
# - Update sheetView selections on a couple of sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "CreateOwn Test Data" (currently Worksheets(1)) gains columns G:O.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 1 - headers (existing header style is already applied on F1, and the
# newly written cells inherit the same style by copying the format across
# before writing the new header labels).
$ws1.Range("F1").Value = "expectedActivityCreatedMsg"
$ws1.Range("F1").Copy()
$ws1.Range("G1:O1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("G1").Value = "activityType"
$ws1.Range("H1").Value = "timeRequired"
$ws1.Range("I1").Value = "date"
$ws1.Range("J1").Value = "completeWithin"
$ws1.Range("K1").Value = "endAfterDays"
$ws1.Range("L1").Value = "expectedAssignedActivityMsg"
$ws1.Range("M1").Value = "points"
$ws1.Range("N1").Value = "statusOfActivity"
$ws1.Range("O1").Value = "expectedToastMsg"

# Row 2 - "Once" activity scenario.
$ws1.Range("G2").Value = "Once"
$ws1.Range("H2").Value = "'30"
$ws1.Range("I2").Value = "'16"
$ws1.Range("J2").Value = "DAY"
$ws1.Range("K2").Value = "'"
$ws1.Range("L2").Value = "YAY! IT'S ASSIGNED"
$ws1.Range("M2").Value = "'5"
$ws1.Range("N2").Value = "Completed"
$ws1.Range("O2").Value = "Activity Deleted Successfully"

# Row 3 - "Every Day" activity scenario.
$ws1.Range("G3").Value = "Every Day"
$ws1.Range("H3").Value = "'60"
$ws1.Range("I3").Value = "'16"
$ws1.Range("J3").Value = "'"
$ws1.Range("K3").Value = "'10"
$ws1.Range("L3").Value = "YAY! IT'S ASSIGNED"
$ws1.Range("M3").Value = "'5"
$ws1.Range("N3").Value = "Completed"
$ws1.Range("O3").Value = "Activity Deleted Successfully"

# Column widths for the newly populated columns.
$ws1.Columns.Item(7).ColumnWidth = 10.77734375
$ws1.Columns.Item(8).ColumnWidth = 11.88671875
$ws1.Columns.Item(10).ColumnWidth = 14.109375
$ws1.Columns.Item(11).ColumnWidth = 14.109375
$ws1.Columns.Item(12).ColumnWidth = 25.88671875
$ws1.Columns.Item(14).ColumnWidth = 14
$ws1.Columns.Item(15).ColumnWidth = 24.21875

# ---------------------------------------------------------------------
# 2) Insert a new sheet ("Sheet1") right after "CreateOwn Test Data" and
#    before "AddAchievement Test Data", holding a single data row that
#    reuses the "Playing Chess" activity values.
# ---------------------------------------------------------------------
$wsAdd = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($wsAdd)

$newSheet.Range("A1").Value = "satish.kshirsagar@gmail.com"
$newSheet.Range("B1").Value = "'1111"
$newSheet.Range("C1").Value = "Playing Chess"
$newSheet.Range("D1").Value = "Playing Chess with friends"
$newSheet.Range("E1").Value = "Chess"
$newSheet.Range("F1").Value = "ACTIVITY CREATED SUCCESSFULLY"
$newSheet.Rows.Item(1).Select()

# ---------------------------------------------------------------------
# 3) sheetView tweaks on the other sheets.
# ---------------------------------------------------------------------
# "Create Activity Error Msg Data" (now the 4th tab) loses tabSelected and
# moves its active cell to F10.
$wsErr = $wb.Worksheets.Item(4)
$wsErr.Range("F10").Select()

# "CreateOwn Test Data" becomes the active tab again with J7 selected.
$ws1.Activate()
$ws1.Range("J7").Select()
